$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24.277779
$ws.Range("I11").Value = 24.277779
$ws.Range("K11").Value = 24.277779
$ws.Range("M11").Value = 115.722221

$ws.Range("H17").Value = 1794.3438
$ws.Range("I17").Value = 299.0909
$ws.Range("J17").Value = 2577.5715
$ws.Range("K17").Value = 897.2727
$ws.Range("L17").Value = 7732.7145
$ws.Range("M17").Value = -729.2727
$ws.Range("N17").Value = -8068.7145

$ws.Range("H28").Value = 20499.8
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20499.8
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 20499.8
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = -21469.8

$ws.Range("H38").Value = 1441.6154
$ws.Range("I38").Value = 360.75
$ws.Range("J38").Value = 1922
$ws.Range("K38").Value = 1082.25
$ws.Range("L38").Value = 5766
$ws.Range("M38").Value = -710.25
$ws.Range("N38").Value = -6510

$ws.Range("H131").Value = 3734.5
$ws.Range("I131").Value = 3734.5
$ws.Range("K131").Value = 11203.5
$ws.Range("M131").Value = -6163.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 359
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 488.5
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 488.5
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -712.5

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null

$ws.Range("H63").Value = 3199.4
$ws.Range("I63").Value = 4499
$ws.Range("J63").Value = 2333
$ws.Range("K63").Value = 4499
$ws.Range("L63").Value = 2333
$ws.Range("M63").Value = -3813
$ws.Range("N63").Value = -3705

$ws.Range("H66").Value = 3199.4
$ws.Range("I66").Value = 4499
$ws.Range("J66").Value = 2333
$ws.Range("K66").Value = 22495
$ws.Range("L66").Value = 11665
$ws.Range("M66").Value = -19063
$ws.Range("N66").Value = -18529

$ws.Range("H122").Value = 1259.125
$ws.Range("I122").Value = 1259.125
$ws.Range("K122").Value = 3777.375
$ws.Range("M122").Value = -1327.375

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 359
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 488.5
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 488.5
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -718.5

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null

$ws.Range("H95").Value = 20312
$ws.Range("J95").Value = 20312
$ws.Range("L95").Value = 20312
$ws.Range("N95").Value = -25804

$ws.Range("H134").Value = 4676.737
$ws.Range("I134").Value = 1777.2667
$ws.Range("K134").Value = 5331.800099999999
$ws.Range("M134").Value = -2796.800099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 807.5
$ws.Range("I5").Value = 135.6
$ws.Range("J5").Value = 2487.25
$ws.Range("K5").Value = 135.6
$ws.Range("L5").Value = 2487.25
$ws.Range("M5").Value = -23.59999999999999
$ws.Range("N5").Value = -2711.25

$ws.Range("H22").Value = 2253.75
$ws.Range("I22").Value = 695
$ws.Range("J22").Value = 3812.5
$ws.Range("K22").Value = 695
$ws.Range("L22").Value = 3812.5
$ws.Range("M22").Value = -345
$ws.Range("N22").Value = -4512.5

$ws.Range("H51").Value = 45588.332
$ws.Range("J51").Value = 45588.332
$ws.Range("L51").Value = 45588.332
$ws.Range("N51").Value = -47060.332

$ws.Range("H59").Value = 38416
$ws.Range("J59").Value = 38416
$ws.Range("L59").Value = 38416
$ws.Range("N59").Value = -40706

$ws.Range("H61").Value = 45588.332
$ws.Range("J61").Value = 45588.332
$ws.Range("L61").Value = 45588.332
$ws.Range("N61").Value = -46284.332

$ws.Range("H99").Value = 4493.25
$ws.Range("I99").Value = 4493.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4493.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2995.25
$ws.Range("N99").Value = $null

$ws.Range("H126").Value = 4493.25
$ws.Range("I126").Value = 4493.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13479.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11009.75
$ws.Range("N126").Value = $null

$ws.Range("H134").Value = 3674.6667
$ws.Range("I134").Value = 2955.7222
$ws.Range("K134").Value = 8867.1666
$ws.Range("M134").Value = -6332.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2619.6
$ws.Range("I12").Value = 2775
$ws.Range("K12").Value = 2775
$ws.Range("M12").Value = -2635

$ws.Range("H113").Value = 6726
$ws.Range("I113").Value = 996.5
$ws.Range("K113").Value = 996.5
$ws.Range("M113").Value = 1173.5

$ws.Range("H122").Value = 1601.7142
$ws.Range("I122").Value = 1822.4
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 5467.200000000001
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -3017.200000000001
$ws.Range("N122").Value = -8050

$ws.Range("H126").Value = 5032.6
$ws.Range("I126").Value = 4353
$ws.Range("J126").Value = 5485.6665
$ws.Range("K126").Value = 13059
$ws.Range("L126").Value = 16456.9995
$ws.Range("M126").Value = -10589
$ws.Range("N126").Value = -21396.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3800.4
$ws.Range("I7").Value = 3250.5
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 3250.5
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -3138.5
$ws.Range("N7").Value = -6224

$ws.Range("H61").Value = 125005976
$ws.Range("I61").Value = 200004770
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 200004770
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -200004568
$ws.Range("N61").Value = -8404

$ws.Range("H68").Value = 7978.5
$ws.Range("J68").Value = 8472.111000000001
$ws.Range("L68").Value = 8472.111000000001
$ws.Range("N68").Value = -9970.111000000001

$ws.Range("H71").Value = 7978.5
$ws.Range("J71").Value = 8472.111000000001
$ws.Range("L71").Value = 42360.55500000001
$ws.Range("N71").Value = -49848.55500000001

$ws.Range("H113").Value = 125005976
$ws.Range("I113").Value = 200004770
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 200004770
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -200002600
$ws.Range("N113").Value = -12340

$ws.Range("H126").Value = 3800.4
$ws.Range("I126").Value = 3250.5
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 9751.5
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -7281.5
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 7865.8667
$ws.Range("I132").Value = 8227.714
$ws.Range("K132").Value = 24683.142
$ws.Range("M132").Value = -22153.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1845
$ws.Range("I13").Value = 517.5
$ws.Range("J13").Value = 4500
$ws.Range("K13").Value = 517.5
$ws.Range("L13").Value = 4500
$ws.Range("M13").Value = -377.5
$ws.Range("N13").Value = -4780

$ws.Range("H51").Value = 21151.25
$ws.Range("I51").Value = 21415.715
$ws.Range("J51").Value = 19300
$ws.Range("K51").Value = 21415.715
$ws.Range("L51").Value = 19300
$ws.Range("M51").Value = -20905.715
$ws.Range("N51").Value = -20320

$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982

$ws.Range("H113").Value = 1292
$ws.Range("I113").Value = 1227.6
$ws.Range("J113").Value = 1399.3334
$ws.Range("K113").Value = 3682.8
$ws.Range("L113").Value = 4198.0002
$ws.Range("M113").Value = -1512.8
$ws.Range("N113").Value = -8538.0002

$ws.Range("H122").Value = 3212.8572
$ws.Range("I122").Value = 2925
$ws.Range("J122").Value = 3596.6667
$ws.Range("K122").Value = 8775
$ws.Range("L122").Value = 10790.0001
$ws.Range("M122").Value = -6325
$ws.Range("N122").Value = -15690.0001

$ws.Range("H132").Value = 1204.1578
$ws.Range("I132").Value = 1076.6666
$ws.Range("K132").Value = 3229.9998
$ws.Range("M132").Value = -699.9998000000001
